$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40 (ALC) - item ID context
$ws.Range("H40").Value = 1997.7407
$ws.Range("I40").Value = 1522.3334
$ws.Range("J40").Value = 2592
$ws.Range("K40").Value = 1522.3334
$ws.Range("L40").Value = 2592
$ws.Range("M40").Value = -1347.3334
$ws.Range("N40").Value = -2942

# Row 76 (ALC) - item ID context
$ws.Range("H76").Value = 151535.67
$ws.Range("I76").Value = 320193.16
$ws.Range("J76").Value = 3960.375
$ws.Range("K76").Value = 320193.16
$ws.Range("L76").Value = 3960.375
$ws.Range("M76").Value = -319878.16
$ws.Range("N76").Value = -4590.375

# Row 79 (ALC) - item ID context
$ws.Range("H79").Value = 151535.67
$ws.Range("I79").Value = 320193.16
$ws.Range("J79").Value = 3960.375
$ws.Range("K79").Value = 320193.16
$ws.Range("L79").Value = 3960.375
$ws.Range("M79").Value = -319101.16
$ws.Range("N79").Value = -6144.375

# Row 107 (ALC) - item ID context
$ws.Range("H107").Value = 255
$ws.Range("I107").Value = 255
$ws.Range("K107").Value = 255
$ws.Range("M107").Value = 1665

# Row 116 (ALC) - item ID context
$ws.Range("H116").Value = 2139563.8
$ws.Range("I116").Value = 2383685.2
$ws.Range("J116").Value = 3500
$ws.Range("K116").Value = 2383685.2
$ws.Range("L116").Value = 3500
$ws.Range("M116").Value = -2380243.2
$ws.Range("N116").Value = -10384

$ws = $wb.Worksheets.Item("ARM")
# Row 63 (ARM) - item ID context
$ws.Range("H63").Value = 2371
$ws.Range("I63").Value = 2356.6667
$ws.Range("K63").Value = 2356.6667
$ws.Range("M63").Value = -1670.6667

# Row 66 (ARM) - item ID context
$ws.Range("H66").Value = 2371
$ws.Range("I66").Value = 2356.6667
$ws.Range("K66").Value = 11783.3335
$ws.Range("M66").Value = -8351.333500000001

# Row 74 (ARM) - item ID context
$ws.Range("H74").Value = 531045
$ws.Range("I74").Value = 5493.4165
$ws.Range("J74").Value = 1431990.6
$ws.Range("K74").Value = 5493.4165
$ws.Range("L74").Value = 1431990.6
$ws.Range("M74").Value = -4619.4165
$ws.Range("N74").Value = -1433738.6

# Row 77 (ARM) - item ID context
$ws.Range("H77").Value = 531045
$ws.Range("I77").Value = 5493.4165
$ws.Range("J77").Value = 1431990.6
$ws.Range("K77").Value = 27467.0825
$ws.Range("L77").Value = 7159953
$ws.Range("M77").Value = -23099.0825
$ws.Range("N77").Value = -7168689

# Row 113 (ARM) - item ID context
$ws.Range("H113").Value = 42000
$ws.Range("J113").Value = 42000
$ws.Range("L113").Value = 42000
$ws.Range("N113").Value = -50678

# Row 132 (ARM) - item ID context
$ws.Range("H132").Value = 19082.22
$ws.Range("I132").Value = 22855.312
$ws.Range("J132").Value = 2617.818
$ws.Range("K132").Value = 68565.936
$ws.Range("L132").Value = 7853.454000000001
$ws.Range("M132").Value = -66035.936
$ws.Range("N132").Value = -12913.454

$ws = $wb.Worksheets.Item("BSM")
# Row 86 (BSM) - item ID context
$ws.Range("H86").Value = 15387270
$ws.Range("I86").Value = 20002480
$ws.Range("J86").Value = 3235.6667
$ws.Range("K86").Value = 20002480
$ws.Range("L86").Value = 3235.6667
$ws.Range("M86").Value = -20001357
$ws.Range("N86").Value = -5481.6667

# Row 89 (BSM) - item ID context
$ws.Range("H89").Value = 15387270
$ws.Range("I89").Value = 20002480
$ws.Range("J89").Value = 3235.6667
$ws.Range("K89").Value = 100012400
$ws.Range("L89").Value = 16178.3335
$ws.Range("M89").Value = -100006784
$ws.Range("N89").Value = -27410.3335

# Row 107 (BSM) - item ID context
$ws.Range("H107").Value = 1049.5
$ws.Range("I107").Value = 1049.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1049.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 870.5
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (CRP) - item ID context
$ws.Range("H16").Value = 2028.9445
$ws.Range("I16").Value = 2150.9
$ws.Range("J16").Value = 1876.5
$ws.Range("K16").Value = 2150.9
$ws.Range("L16").Value = 1876.5
$ws.Range("M16").Value = -1863.9
$ws.Range("N16").Value = -2450.5

# Row 31 (CRP) - item ID context
$ws.Range("H31").Value = 2619.4736
$ws.Range("I31").Value = 1470.1428
$ws.Range("J31").Value = 3289.9167
$ws.Range("K31").Value = 1470.1428
$ws.Range("L31").Value = 3289.9167
$ws.Range("M31").Value = -1175.1428
$ws.Range("N31").Value = -3879.9167

# Row 34 (CRP) - item ID context
$ws.Range("H34").Value = 2619.4736
$ws.Range("I34").Value = 1470.1428
$ws.Range("J34").Value = 3289.9167
$ws.Range("K34").Value = 1470.1428
$ws.Range("L34").Value = 3289.9167
$ws.Range("M34").Value = -1268.1428
$ws.Range("N34").Value = -3693.9167

# Row 58 (CRP) - item ID context
$ws.Range("H58").Value = 2419.9744
$ws.Range("I58").Value = 880.619
$ws.Range("J58").Value = 4215.8887
$ws.Range("K58").Value = 880.619
$ws.Range("L58").Value = 4215.8887
$ws.Range("M58").Value = -677.619
$ws.Range("N58").Value = -4621.8887

# Row 107 (CRP) - item ID context
$ws.Range("H107").Value = 1208.1177
$ws.Range("I107").Value = 1146.1
$ws.Range("J107").Value = 1296.7142
$ws.Range("K107").Value = 1146.1
$ws.Range("L107").Value = 1296.7142
$ws.Range("M107").Value = 773.9000000000001
$ws.Range("N107").Value = -5136.7142

# Row 113 (CRP) - item ID context
$ws.Range("H113").Value = 2028.9445
$ws.Range("I113").Value = 2150.9
$ws.Range("J113").Value = 1876.5
$ws.Range("K113").Value = 2150.9
$ws.Range("L113").Value = 1876.5
$ws.Range("M113").Value = 19.09999999999991
$ws.Range("N113").Value = -6216.5

# Row 136 (CRP) - item ID context
$ws.Range("H136").Value = 2419.9744
$ws.Range("I136").Value = 880.619
$ws.Range("J136").Value = 4215.8887
$ws.Range("K136").Value = 2641.857
$ws.Range("L136").Value = 12647.6661
$ws.Range("M136").Value = -91.85699999999997
$ws.Range("N136").Value = -17747.6661

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (CUL) - item ID context
$ws.Range("H5").Value = 839.75
$ws.Range("I5").Value = 503.73914
$ws.Range("J5").Value = 1434.2307
$ws.Range("K5").Value = 1511.21742
$ws.Range("L5").Value = 4302.6921
$ws.Range("M5").Value = -1399.21742
$ws.Range("N5").Value = -4526.6921

# Row 107 (CUL) - item ID context
$ws.Range("H107").Value = 88.96154
$ws.Range("I107").Value = 95.30768999999999
$ws.Range("J107").Value = 82.61539
$ws.Range("K107").Value = 285.92307
$ws.Range("L107").Value = 247.84617
$ws.Range("M107").Value = 1634.07693
$ws.Range("N107").Value = -4087.84617

# Row 132 (CUL) - item ID context
$ws.Range("H132").Value = 1292.2916
$ws.Range("J132").Value = 1763.8
$ws.Range("L132").Value = 15874.2
$ws.Range("N132").Value = -20934.2

# Row 134 (CUL) - item ID context
$ws.Range("H134").Value = 844.3043
$ws.Range("I134").Value = 551.3570999999999
$ws.Range("K134").Value = 1654.0713
$ws.Range("M134").Value = 3415.9287

# Row 135 (CUL) - item ID context
$ws.Range("H135").Value = 839.75
$ws.Range("I135").Value = 503.73914
$ws.Range("J135").Value = 1434.2307
$ws.Range("K135").Value = 4533.65226
$ws.Range("L135").Value = 12908.0763
$ws.Range("M135").Value = -1998.65226
$ws.Range("N135").Value = -17978.0763

# Row 137 (CUL) - item ID context
$ws.Range("H137").Value = 2902.606
$ws.Range("I137").Value = 811.3889
$ws.Range("J137").Value = 5412.067
$ws.Range("K137").Value = 2434.1667
$ws.Range("L137").Value = 16236.201
$ws.Range("M137").Value = 2665.8333
$ws.Range("N137").Value = -26436.201

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (GSM) - item ID context
$ws.Range("H70").Value = 2107247.2
$ws.Range("I70").Value = 4205923.5
$ws.Range("J70").Value = 8571.429
$ws.Range("K70").Value = 4205923.5
$ws.Range("L70").Value = 8571.429
$ws.Range("M70").Value = -4205653.5
$ws.Range("N70").Value = -9111.429

# Row 73 (GSM) - item ID context
$ws.Range("H73").Value = 2107247.2
$ws.Range("I73").Value = 4205923.5
$ws.Range("J73").Value = 8571.429
$ws.Range("K73").Value = 4205923.5
$ws.Range("L73").Value = 8571.429
$ws.Range("M73").Value = -4204987.5
$ws.Range("N73").Value = -10443.429

# Row 80 (GSM) - item ID context
$ws.Range("H80").Value = 4051.7856
$ws.Range("I80").Value = 6411
$ws.Range("K80").Value = 6411
$ws.Range("M80").Value = -5413

# Row 83 (GSM) - item ID context
$ws.Range("H83").Value = 4051.7856
$ws.Range("I83").Value = 6411
$ws.Range("K83").Value = 32055
$ws.Range("M83").Value = -27063

# Row 113 (GSM) - item ID context
$ws.Range("H113").Value = 2333.3333
$ws.Range("I113").Value = 2250
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 2250
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = -80
$ws.Range("N113").Value = -6840

$ws = $wb.Worksheets.Item("LTW")
# Row 19 (LTW) - item ID context
$ws.Range("H19").Value = 2394.3333
$ws.Range("I19").Value = 1091.5
$ws.Range("K19").Value = 1091.5
$ws.Range("M19").Value = -921.5

# Row 55 (LTW) - item ID context
$ws.Range("H55").Value = 252.85294
$ws.Range("I55").Value = 310.89474
$ws.Range("J55").Value = 179.33333
$ws.Range("K55").Value = 310.89474
$ws.Range("L55").Value = 179.33333
$ws.Range("M55").Value = -137.89474
$ws.Range("N55").Value = -525.3333299999999

# Row 93 (LTW) - item ID context
$ws.Range("H93").Value = 1095.6774
$ws.Range("I93").Value = 1079.8948
$ws.Range("J93").Value = 1120.6666
$ws.Range("K93").Value = 1079.8948
$ws.Range("L93").Value = 1120.6666
$ws.Range("M93").Value = 168.1052
$ws.Range("N93").Value = -3616.6666

$ws = $wb.Worksheets.Item("WVR")
# Row 31 (WVR) - item ID context
$ws.Range("H31").Value = 4001
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 4001
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 4001
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -4697

# Row 113 (WVR) - item ID context
$ws.Range("H113").Value = 18518798
$ws.Range("I113").Value = 266.26086
$ws.Range("J113").Value = 125000350
$ws.Range("K113").Value = 798.7825799999999
$ws.Range("L113").Value = 375001050
$ws.Range("M113").Value = 1371.21742
$ws.Range("N113").Value = -375005390

